# data updates / fixes
# - Oakland Raiders are relocating to Las Vegas: update the existing
#   "Oakland Raiders" row's pfr/fff abbreviation codes (B26/C26) to the
#   (placeholder) "LVE" / "xxx" values, and append a brand-new row for the
#   "Las Vegas Raiders" franchise entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26 (Oakland Raiders): update pfr_abbreviation / fff_abbreviation ---
$ws.Range("B26").Value = "LVE"
$ws.Range("C26").Value = "xxx"

# --- New row 34: Las Vegas Raiders ---
$ws.Range("A34").Value = "Las Vegas Raiders"
$ws.Range("B34").Value = "LVE"
$ws.Range("C34").Value = "OAK"
$ws.Range("D34").Value = "Las_Vegas"

# Match the formatting of the row above it (thin right border, vertical
# centered, no bottom border) rather than leaving the new row unstyled.
$ws.Range("D33").Copy() | Out-Null
$ws.Range("D34").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Reflect the edit location in the view, same as the author's saved selection.
$ws.Range("B27").Select() | Out-Null
